# Update the division-problem worksheet: replace each problem's text with
# the new value from the commit, cell by cell, preserving cell formatting.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "50÷6="
$t.Cell(1, 2).Range.Text = "54÷4="
$t.Cell(1, 3).Range.Text = "92÷6="
$t.Cell(1, 4).Range.Text = "32÷8="
$t.Cell(1, 5).Range.Text = "50÷2="

$t.Cell(5, 1).Range.Text = "76÷2="
$t.Cell(5, 2).Range.Text = "11÷6="
$t.Cell(5, 3).Range.Text = "75÷4="
$t.Cell(5, 4).Range.Text = "28÷9="
$t.Cell(5, 5).Range.Text = "15÷4="

$t.Cell(9, 1).Range.Text = "41÷5="
$t.Cell(9, 2).Range.Text = "74÷3="
$t.Cell(9, 3).Range.Text = "92÷8="
$t.Cell(9, 4).Range.Text = "56÷2="
$t.Cell(9, 5).Range.Text = "79÷2="

$t.Cell(13, 1).Range.Text = "19÷9="
$t.Cell(13, 2).Range.Text = "62÷8="
$t.Cell(13, 3).Range.Text = "82÷9="
$t.Cell(13, 4).Range.Text = "33÷6="
$t.Cell(13, 5).Range.Text = "84÷4="

$t.Cell(17, 1).Range.Text = "98÷2="
$t.Cell(17, 2).Range.Text = "27÷4="
$t.Cell(17, 3).Range.Text = "51÷9="
$t.Cell(17, 4).Range.Text = "50÷4="
$t.Cell(17, 5).Range.Text = "71÷3="
